# Grade tables now also in stripped form for easier input.
# Fill column A (the "X" marker column) for rows 10 through 50 with "X",
# matching the value already present in row 9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 10; $r -le 50; $r++) {
    $ws.Cells.Item($r, 1).Value = "X"
}

# Reset the active cell selection to A1 (matches the saved view state).
$ws.Range("A1").Select()
